# Reconstruction bdd avec page edition
# Update the row-2 data values (commune statistics) on the active worksheet
# to reflect the refreshed dataset, and drop the now-removed BC2 entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 114953
$ws.Range("B2").Value = 22825
$ws.Range("C2").Value = 19330
$ws.Range("D2").Value = 99355
$ws.Range("E2").Value = 24546
$ws.Range("J2").Value = 59041
$ws.Range("K2").Value = 75250
$ws.Range("L2").Value = 17408
$ws.Range("M2").Value = 105144
$ws.Range("P2").Value = 107401
$ws.Range("S2").Value = 14770
$ws.Range("T2").Value = 90649
$ws.Range("V2").Value = 311343
$ws.Range("W2").Value = 19914
$ws.Range("Z2").Value = 51079
$ws.Range("AA2").Value = 62734
$ws.Range("AC2").Value = 57215
$ws.Range("AD2").Value = 190793
$ws.Range("AE2").Value = 72022
$ws.Range("AF2").Value = 84524
$ws.Range("AG2").Value = 54576
$ws.Range("AH2").Value = 54712
$ws.Range("AI2").Value = 83137
$ws.Range("AJ2").Value = 27943
$ws.Range("AK2").Value = 46099
$ws.Range("AL2").Value = 265414
$ws.Range("AP2").Value = 228532
$ws.Range("AQ2").Value = 19242
$ws.Range("AR2").Value = 578330
$ws.Range("AT2").Value = 521829
$ws.Range("AU2").Value = 13085
$ws.Range("AW2").Value = 16327
$ws.Range("AX2").Value = 28783
$ws.Range("AY2").Value = 93448
$ws.Range("AZ2").Value = 58382
$ws.Range("BC2").ClearContents()
$ws.Range("BD2").Value = 42817
$ws.Range("BE2").Value = 214506
$ws.Range("BF2").Value = 217250
$ws.Range("BG2").Value = 55385
$ws.Range("BH2").Value = 27839
$ws.Range("BI2").Value = 104840
$ws.Range("BK2").Value = 53300
$ws.Range("BL2").Value = 18173
$ws.Range("BO2").Value = 138632
$ws.Range("BP2").Value = 179069
$ws.Range("BQ2").Value = 54860
$ws.Range("BR2").Value = 54249
$ws.Range("BS2").Value = 75079
$ws.Range("BT2").Value = 24324
$ws.Range("BU2").Value = 57702
$ws.Range("BV2").Value = 113788
$ws.Range("BZ2").Value = 24807
$ws.Range("CE2").Value = 22399
$ws.Range("CG2").Value = 72883
$ws.Range("CH2").Value = 19794
$ws.Range("CI2").Value = 115591
$ws.Range("CJ2").Value = 128787
$ws.Range("CK2").Value = 136464
$ws.Range("CL2").Value = 20509
$ws.Range("CM2").Value = 24542
$ws.Range("CN2").Value = 136528
$ws.Range("CQ2").Value = 93679
$ws.Range("CR2").Value = 125623
$ws.Range("CS2").Value = 84677
$ws.Range("CT2").Value = 27442
$ws.Range("CU2").Value = 137889
$ws.Range("CV2").Value = 199897
$ws.Range("CW2").Value = 138701
$ws.Range("CX2").Value = 142391
$ws.Range("DB2").Value = 54680
$ws.Range("DC2").Value = 62957
$ws.Range("DE2").Value = 46387
